$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.468.96"

$ws.Range("D3").Value = "3.844.65"
$ws.Range("E3").Value = "  -1.26%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.24%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.60%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.79%  "

$ws.Range("D7").Value = "3.844.07"
$ws.Range("E7").Value = "  -1.33%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("E9").Value = "  -1.40%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.167"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.32%  "

$ws.Range("E11").Value = "  +1.57%  "

$ws.Range("E12").Value = "  -1.86%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000269"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.10%  "

$ws.Range("E14").Value = "  -2.83%  "

$ws.Range("D15").Value = "4.490.11"
$ws.Range("E15").Value = "  -1.25%  "

$ws.Range("D16").Value = "3.830.16"
$ws.Range("E16").Value = "  -1.28%  "

$ws.Range("D17").Value = "68.453.27"
$ws.Range("E17").Value = "  -1.79%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.51"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.42%  "

$ws.Range("E19").Value = "  -3.32%  "

$ws.Range("E20").Value = "  -0.86%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.08"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.07%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "470.67"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.06%  "

$ws.Range("E23").Value = "  -1.80%  "

$ws.Range("E24").Value = "  -3.66%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.33"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.33%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.23"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.61%  "

$ws.Range("E27").Value = "  -1.39%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.07"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.78%  "

$ws.Range("E29").Value = "  +0.05%  "

$ws.Range("E30").Value = "  -0.28%  "

$ws.Range("D31").Value = "3.991.14"
$ws.Range("E31").Value = "  -1.25%  "

$ws.Range("E32").Value = "  -2.16%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.65"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.95%  "

$ws.Range("E34").Value = "  -4.67%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.39"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.51%  "

$ws.Range("D36").Value = "3.806.33"
$ws.Range("E36").Value = "  -1.34%  "

$ws.Range("E37").Value = "  -2.15%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.70"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +10.37%  "

$ws.Range("E39").Value = "  -1.12%  "

$ws.Range("E40").Value = "  -1.33%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.96"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.87%  "

$ws.Range("E42").Value = "  +0.07%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.315"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.24%  "

$ws.Range("E44").Value = "  -5.23%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.74"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.52%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "415.64"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.27%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "47.16"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.32%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000290"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.96%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0360"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.12%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "141.65"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.76%  "
